# RSTK-9619-Derived Component Receipt Reversal.xlsx — "updated SYDATA-Work order testcases"
#
# 1) "Create Disassembly WO": insert a new "Site" column at A, reorder/update the
#    two disassembly-item rows (qty + flags), re-fit column widths.
# 2) "Derived Components": stamp a "Disassy Good Loc ID" value ("OH (On Hand Loc)")
#    in column F for the 4 data rows (small Consolas font), and bump a couple of
#    qty values.
# 3) "Consumable Components": add explicit 0 Scrap Factor / Setup Qty values on the
#    two "MS-serial and lot tracked" rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Create Disassembly WO"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Create Disassembly WO")

# Insert a new first column ("Site") — everything else shifts right by one.
$ws1.Columns.Item(1).Insert()

$ws1.Range("A1").Value = "Site"
$ws1.Range("A2").Value = "10 (Denver)"
$ws1.Range("A3").Value = "10 (Denver)"

# Re-key the two data rows to their final (post-edit) content. Row 2 now holds
# the "Serial" disassembly item (qty bumped 2 -> 3, "Add derived" flipped on);
# row 3 holds the "Lot and serial" item unchanged apart from the shift.
$ws1.Range("B2").Value = "Pro-Disassembley Serial (Serial track)"
$ws1.Range("C2").Value = 3
$ws1.Range("D2").Value = "100 Home Project"
$ws1.Range("F2").Value = $true
$ws1.Range("G2").Value = $true

$ws1.Range("B3").Value = "Pro-Disassembley (Lot and serial track)"
$ws1.Range("C3").Value = 3
$ws1.Range("D3").Value = "100 Home Project"
$ws1.Range("F3").Value = $true
$ws1.Range("G3").Value = $false

# Column widths (best-fit re-computation after the insert/content changes).
$ws1.Columns.Item(1).ColumnWidth = 9.498697916666666
$ws1.Columns.Item(4).ColumnWidth = 14.830729166666666
$ws1.Columns.Item(5).ColumnWidth = 16.944010416666668
$ws1.Columns.Item(7).ColumnWidth = 25.385416666666668

$ws1.Range("B16").Select()

# ---------------------------------------------------------------------------
# Sheet 2: "Derived Components"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Derived Components")

$locCell = $ws2.Range("F2")
$locCell.Value = "OH (On Hand Loc)"
$locCell.Font.Size = 7
$locCell.Font.Color = 2367776
$locCell.Font.Name = "Consolas"
$locCell.Font.Family = 3

$ws2.Range("F3").Value = "OH (On Hand Loc)"
$ws2.Range("F4").Value = "OH (On Hand Loc)"
$ws2.Range("F5").Value = "OH (On Hand Loc)"

$locCell.Copy()
$ws2.Range("F3:F5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("C4").Value = 3
$ws2.Range("D4").Value = 3
$ws2.Range("C5").Value = 3
$ws2.Range("D5").Value = 3

$ws2.Range("G17").Select()

# ---------------------------------------------------------------------------
# Sheet 3: "Consumable Components"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Consumable Components")

$ws3.Range("G2").Value = 0
$ws3.Range("H2").Value = 0
$ws3.Range("G5").Value = 0
$ws3.Range("H5").Value = 0

$ws3.Range("C15").Select()

# Restore sheet 1 as the active tab/selection.
$ws1.Activate()
$ws1.Range("B16").Select()
